# Add Buyung (B. Hadi) to the list of authors, appending after "Sander".
#
# The author line currently reads:
#   A.H. Sparks and N.P. Castilla and B.O. Sander
# and should become:
#   A.H. Sparks and N.P. Castilla and B.O. Sander and B. Hadi
#
# Each "word"/space in that paragraph lives in its own run (that's how the
# rest of the author paragraph is structured), so we insert the new pieces
# one at a time right after "Sander" to keep the same one-run-per-token
# pattern instead of dumping everything into a single new run.

$d = $word.ActiveDocument

$range = $d.Content
$found = $range.Find.Execute("Sander", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)

if ($found) {
    # Collapse the found range to its end (right after "Sander") and append
    # the new tokens, each as its own run, in order.
    $range.Collapse(0)

    $range.InsertAfter(" ")
    $range.Collapse(0)

    $range.InsertAfter("and")
    $range.Collapse(0)

    $range.InsertAfter(" ")
    $range.Collapse(0)

    $range.InsertAfter("B.")
    $range.Collapse(0)

    $range.InsertAfter(" ")
    $range.Collapse(0)

    $range.InsertAfter("Hadi")
    $range.Collapse(0)
}
